# Re-run of the CircaDiPy cosinor-per-day analysis (sine_0.1, fixed period 1).
# Refreshes the per-day cosinor statistics (columns F..Z) and the corresponding
# "significant" flags (column B) for rows 2-13 to match the latest simulation run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0
$ws.Range("F2").Value = 22.77000000000012
$ws.Range("H2").Value = 0.1278635932895194
$ws.Range("I2").Value = 0.1278635932895194
$ws.Range("L2").Value = 24.34961020373348
$ws.Range("M2").Value = "[-4.1750986219058674, 52.87431902937283]"
$ws.Range("N2").Value = 0.09243168241607336
$ws.Range("O2").Value = 0.09243168241607336
$ws.Range("P2").Value = 1.578658170272348
$ws.Range("Q2").Value = "[0.044026323473730145, 3.113290017070966]"
$ws.Range("R2").Value = 0.04403655663985329
$ws.Range("S2").Value = 0.04403655663985329
$ws.Range("T2").Value = 69.67338827129973
$ws.Range("U2").Value = "[53.34317252712222, 86.00360401547724]"
$ws.Range("V2").Value = [double]"4.798472730271897e-11"
$ws.Range("W2").Value = [double]"4.798472730271897e-11"
$ws.Range("X2").Value = 17.0490090090091
$ws.Range("Y2").Value = 11.48756756756763
$ws.Range("Z2").Value = 22.61045045045057

# Row 3
$ws.Range("F3").Value = 22.77000000000012
$ws.Range("H3").Value = [double]"8.281595349046711e-05"
$ws.Range("I3").Value = [double]"8.281595349046711e-05"
$ws.Range("L3").Value = 45.25451309703597
$ws.Range("M3").Value = "[21.090280276131836, 69.4187459179401]"
$ws.Range("N3").Value = 0.0004704183200632261
$ws.Range("O3").Value = 0.0004704183200632261
$ws.Range("P3").Value = 1.880552959806502
$ws.Range("Q3").Value = "[1.2893423303021168, 2.4717635893108882]"
$ws.Range("R3").Value = [double]"7.764234122298319e-08"
$ws.Range("S3").Value = [double]"7.764234122298319e-08"
$ws.Range("T3").Value = 61.53764235988509
$ws.Range("U3").Value = "[48.40896710459118, 74.666317615179]"
$ws.Range("V3").Value = [double]"3.057554209817681e-12"
$ws.Range("W3").Value = [double]"3.057554209817681e-12"
$ws.Range("X3").Value = 15.95495495495504
$ws.Range("Y3").Value = 13.81243243243251
$ws.Range("Z3").Value = 18.09747747747758

# Row 4
$ws.Range("F4").Value = 22.77000000000012
$ws.Range("H4").Value = 0.0002291611244584768
$ws.Range("I4").Value = 0.0002291611244584768
$ws.Range("L4").Value = 35.88436904331667
$ws.Range("M4").Value = "[14.193880644646505, 57.57485744198683]"
$ws.Range("N4").Value = 0.001730013369199757
$ws.Range("O4").Value = 0.001730013369199757
$ws.Range("P4").Value = 2.182447749340657
$ws.Range("Q4").Value = "[1.6038160694001942, 2.7610794292811196]"
$ws.Range("R4").Value = [double]"1.340329847110411e-09"
$ws.Range("S4").Value = [double]"1.340329847110411e-09"
$ws.Range("T4").Value = 57.25878047146323
$ws.Range("U4").Value = "[45.962914327062094, 68.55464661586437]"
$ws.Range("V4").Value = [double]"2.704503287986881e-13"
$ws.Range("W4").Value = [double]"2.704503287986881e-13"
$ws.Range("X4").Value = 14.86090090090098
$ws.Range("Y4").Value = 12.76396396396403
$ws.Range("Z4").Value = 16.95783783783793

# Row 5
$ws.Range("F5").Value = 22.77000000000012
$ws.Range("H5").Value = [double]"6.287037146446295e-08"
$ws.Range("I5").Value = [double]"6.287037146446295e-08"
$ws.Range("L5").Value = 50.92142183958948
$ws.Range("M5").Value = "[31.804263029172205, 70.03858065000675]"
$ws.Range("N5").Value = [double]"2.705914632583628e-06"
$ws.Range("O5").Value = [double]"2.705914632583628e-06"
$ws.Range("P5").Value = 2.735921530153273
$ws.Range("Q5").Value = "[2.333395144107734, 3.138447916198812]"
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 62.9815242570908
$ws.Range("U5").Value = "[52.4207972134724, 73.54225130070921]"
$ws.Range("V5").Value = [double]"1.110223024625157e-15"
$ws.Range("W5").Value = [double]"1.110223024625157e-15"
$ws.Range("X5").Value = 12.8551351351352
$ws.Range("Y5").Value = 11.39639639639646
$ws.Range("Z5").Value = 14.31387387387395

# Row 6
$ws.Range("F6").Value = 22.77000000000012
$ws.Range("H6").Value = [double]"2.154981729729499e-07"
$ws.Range("I6").Value = [double]"2.154981729729499e-07"
$ws.Range("L6").Value = 70.47220387195722
$ws.Range("M6").Value = "[43.23686217664424, 97.7075455672702]"
$ws.Range("N6").Value = [double]"4.5340508187941e-06"
$ws.Range("O6").Value = [double]"4.5340508187941e-06"
$ws.Range("P6").Value = 2.786237328408966
$ws.Range("Q6").Value = "[2.333395144107734, 3.2390795127101972]"
$ws.Range("R6").Value = [double]"4.440892098500626e-16"
$ws.Range("S6").Value = [double]"4.440892098500626e-16"
$ws.Range("T6").Value = 73.45674315659457
$ws.Range("U6").Value = "[58.00601967915949, 88.90746663402965]"
$ws.Range("V6").Value = [double]"1.987521258683955e-12"
$ws.Range("W6").Value = [double]"1.987521258683955e-12"
$ws.Range("X6").Value = 12.67279279279286
$ws.Range("Y6").Value = 11.03171171171177
$ws.Range("Z6").Value = 14.31387387387395

# Row 7
$ws.Range("F7").Value = 22.77000000000012
$ws.Range("H7").Value = [double]"2.72640369169963e-08"
$ws.Range("I7").Value = [double]"2.72640369169963e-08"
$ws.Range("L7").Value = 62.39658663123243
$ws.Range("M7").Value = "[43.186724702749544, 81.60644855971532]"
$ws.Range("N7").Value = [double]"4.880186699196543e-08"
$ws.Range("O7").Value = [double]"4.880186699196543e-08"
$ws.Range("P7").Value = -3.069263693597235
$ws.Range("Q7").Value = "[-3.459211130078852, -2.679316257115619]"
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 67.13784211170579
$ws.Range("U7").Value = "[54.60930004679889, 79.6663841766127]"
$ws.Range("V7").Value = [double]"4.507505479978136e-14"
$ws.Range("W7").Value = [double]"4.507505479978136e-14"
$ws.Range("X7").Value = 11.12288288288294
$ws.Range("Y7").Value = 9.70972972972978
$ws.Range("Z7").Value = 12.5360360360361

# Row 8
$ws.Range("F8").Value = 25.56000000000056
$ws.Range("H8").Value = 0.001761210413826042
$ws.Range("I8").Value = 0.001761210413826042
$ws.Range("L8").Value = 43.17187732824567
$ws.Range("M8").Value = "[17.39951868242308, 68.94423597406826]"
$ws.Range("N8").Value = 0.001533653867649587
$ws.Range("O8").Value = 0.001533653867649587
$ws.Range("P8").Value = -3.132158441416851
$ws.Range("Q8").Value = "[-3.9120533143800826, -2.3522635684536186]"
$ws.Range("R8").Value = [double]"2.559914502597849e-10"
$ws.Range("S8").Value = [double]"2.559914502597849e-10"
$ws.Range("T8").Value = 57.72845560448977
$ws.Range("U8").Value = "[42.123591635691405, 73.33331957328812]"
$ws.Range("V8").Value = [double]"2.195576387364895e-09"
$ws.Range("W8").Value = [double]"2.195576387364895e-09"
$ws.Range("X8").Value = 12.7416216216219
$ws.Range("Y8").Value = 9.569009009009218
$ws.Range("Z8").Value = 15.91423423423458

# Row 9
$ws.Range("F9").Value = 25.56000000000056
$ws.Range("H9").Value = 0.005982681524594913
$ws.Range("I9").Value = 0.005982681524594913
$ws.Range("L9").Value = 39.02240000755313
$ws.Range("M9").Value = "[11.978951091348023, 66.06584892375824]"
$ws.Range("N9").Value = 0.005657151655191539
$ws.Range("O9").Value = 0.005657151655191539
$ws.Range("P9").Value = 2.987500521431735
$ws.Range("Q9").Value = "[2.0818161528292727, 3.893184890034197]"
$ws.Range("R9").Value = [double]"3.445776775024001e-08"
$ws.Range("S9").Value = [double]"3.445776775024001e-08"
$ws.Range("T9").Value = 72.8564484528124
$ws.Range("U9").Value = "[56.81445309034568, 88.89844381527912]"
$ws.Range("V9").Value = [double]"7.859268791321483e-12"
$ws.Range("W9").Value = [double]"7.859268791321483e-12"
$ws.Range("X9").Value = 13.40684684684714
$ws.Range("Y9").Value = 9.722522522522738
$ws.Range("Z9").Value = 17.09117117117154

# Row 10
$ws.Range("B10").Value = 1
$ws.Range("F10").Value = 25.56000000000056
$ws.Range("H10").Value = 0.000472701856951141
$ws.Range("I10").Value = 0.000472701856951141
$ws.Range("L10").Value = 41.36945758763776
$ws.Range("M10").Value = "[17.931716449993004, 64.80719872528252]"
$ws.Range("N10").Value = 0.0009022329987280475
$ws.Range("O10").Value = 0.0009022329987280475
$ws.Range("P10").Value = 2.949763672739966
$ws.Range("Q10").Value = "[2.257921446724195, 3.6416058987557367]"
$ws.Range("R10").Value = [double]"4.891576033116962e-11"
$ws.Range("S10").Value = [double]"4.891576033116962e-11"
$ws.Range("T10").Value = 64.45444068558
$ws.Range("U10").Value = "[50.93720109182725, 77.97168027933274]"
$ws.Range("V10").Value = [double]"1.816324868286756e-12"
$ws.Range("W10").Value = [double]"1.816324868286756e-12"
$ws.Range("X10").Value = 13.56036036036065
$ws.Range("Y10").Value = 10.74594594594618
$ws.Range("Z10").Value = 16.37477477477513

# Row 11
$ws.Range("B11").Value = 1
$ws.Range("F11").Value = 25.56000000000056
$ws.Range("H11").Value = 0.0003334298603493435
$ws.Range("I11").Value = 0.0003334298603493435
$ws.Range("L11").Value = 48.99499542121059
$ws.Range("M11").Value = "[18.829779987617982, 79.1602108548032]"
$ws.Range("N11").Value = 0.002058626401258756
$ws.Range("O11").Value = 0.002058626401258756
$ws.Range("P11").Value = 2.006342455445733
$ws.Range("Q11").Value = "[1.364816027685654, 2.6478688832058115]"
$ws.Range("R11").Value = [double]"1.122266188868792e-07"
$ws.Range("S11").Value = [double]"1.122266188868792e-07"
$ws.Range("T11").Value = 63.5261810331336
$ws.Range("U11").Value = "[47.092308480385455, 79.96005358588175]"
$ws.Range("V11").Value = [double]"7.082905373323456e-10"
$ws.Range("W11").Value = [double]"7.082905373323456e-10"
$ws.Range("X11").Value = 17.39819819819858
$ws.Range("Y11").Value = 14.78846846846879
$ws.Range("Z11").Value = 20.00792792792837

# Row 12
$ws.Range("F12").Value = 25.56000000000056
$ws.Range("H12").Value = [double]"2.525523542518116e-07"
$ws.Range("I12").Value = [double]"2.525523542518116e-07"
$ws.Range("L12").Value = 57.73072104686539
$ws.Range("M12").Value = "[34.94774549528731, 80.51369659844347]"
$ws.Range("N12").Value = [double]"6.509821331990295e-06"
$ws.Range("O12").Value = [double]"6.509821331990295e-06"
$ws.Range("P12").Value = 1.930868758062195
$ws.Range("Q12").Value = "[1.490605523324886, 2.371131992799503]"
$ws.Range("R12").Value = [double]"2.182232172742715e-11"
$ws.Range("S12").Value = [double]"2.182232172742715e-11"
$ws.Range("T12").Value = 68.91197901254603
$ws.Range("U12").Value = "[56.01183997330202, 81.81211805179004]"
$ws.Range("V12").Value = [double]"4.996003610813204e-14"
$ws.Range("W12").Value = [double]"4.996003610813204e-14"
$ws.Range("X12").Value = 17.70522522522561
$ws.Range("Y12").Value = 15.91423423423458
$ws.Range("Z12").Value = 19.49621621621664

# Row 13
$ws.Range("F13").Value = 25.56000000000056
$ws.Range("H13").Value = 0.0007203310921424233
$ws.Range("I13").Value = 0.0007203310921424233
$ws.Range("L13").Value = 40.87769639857494
$ws.Range("M13").Value = "[16.314751323880515, 65.44064147326937]"
$ws.Range("N13").Value = 0.001634246585420662
$ws.Range("O13").Value = 0.001634246585420662
$ws.Range("P13").Value = 1.30192127986604
$ws.Range("Q13").Value = "[0.5723422051585008, 2.031500354573579]"
$ws.Range("R13").Value = 0.0008033516613883407
$ws.Range("S13").Value = 0.0008033516613883407
$ws.Range("T13").Value = 52.4726512151335
$ws.Range("U13").Value = "[38.1147314974364, 66.8305709328306]"
$ws.Range("V13").Value = [double]"2.98183766567206e-09"
$ws.Range("W13").Value = [double]"2.98183766567206e-09"
$ws.Range("X13").Value = 20.26378378378423
$ws.Range("Y13").Value = 17.29585585585624
$ws.Range("Z13").Value = 23.23171171171222
